$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 220 (D and F values revised) ---
$ws.Range("D220").Value = 4.58647
$ws.Range("F220").Value = 4.5542

# --- Prepare formatting for the new rows by copying row 220's style ---
$ws.Range("A220:G220").Copy()
$ws.Range("A221:G221").PasteSpecial(-4122)
$ws.Range("A220:G220").Copy()
$ws.Range("A222:G222").PasteSpecial(-4122)
$ws.Range("A220:G220").Copy()
$ws.Range("A223:G223").PasteSpecial(-4122)

# --- Row 221 ---
$ws.Range("A221").Value = 45170.33333333334
$ws.Range("B221").Value = "FX_IDC:USDRON"
$ws.Range("C221").Value = 4.5546
$ws.Range("D221").Value = 4.7417
$ws.Range("E221").Value = 4.54653
$ws.Range("F221").Value = 4.7021
$ws.Range("G221").Value = 0

# --- Row 222 ---
$ws.Range("A222").Value = 45201.375
$ws.Range("B222").Value = "FX_IDC:USDRON"
$ws.Range("C222").Value = 4.6964
$ws.Range("D222").Value = 4.7602
$ws.Range("E222").Value = 4.64967
$ws.Range("F222").Value = 4.6949
$ws.Range("G222").Value = 0

# --- Row 223 ---
$ws.Range("A223").Value = 45231.375
$ws.Range("B223").Value = "FX_IDC:USDRON"
$ws.Range("C223").Value = 4.6941
$ws.Range("D223").Value = 4.7224
$ws.Range("E223").Value = 4.61965
$ws.Range("F223").Value = 4.6438
$ws.Range("G223").Value = 0
